$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Quantity for the Keyboard sale on 2025-01-03 (row 4) from 7 to 1
$ws.Range("D4").Value = 1
